# 自动更新Excel文件 - daily refresh of the remaining-day tracker.
#
# "Today" advanced from 2025-12-07 to 2025-12-08 (one day later).
# Column E ("剩余" / remaining days) counts down by 1 per day for every
# tracked row. When a row's remaining-day counter would reach 0, the item
# is treated as replenished/restocked instead: E is reset back to the
# row's total-day allotment (column D, "总天") and F ("开始时间" / start
# date) is reset to the new "today" (2025-12-08).
#
# Row 36 carries a malformed start date (202510929, not a valid yyyymmdd
# value) in the source data, so it cannot be aged and is left untouched,
# exactly as in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; E=7; F=$null},
    @{Row=3; E=7; F=$null},
    @{Row=4; E=7; F=$null},
    @{Row=5; E=7; F=$null},
    @{Row=6; E=7; F=$null},
    @{Row=7; E=7; F=$null},
    @{Row=8; E=7; F=$null},
    @{Row=9; E=7; F=$null},
    @{Row=10; E=7; F=20251208},
    @{Row=11; E=7; F=$null},
    @{Row=12; E=7; F=$null},
    @{Row=13; E=7; F=$null},
    @{Row=14; E=7; F=$null},
    @{Row=15; E=7; F=$null},
    @{Row=16; E=1; F=$null},
    @{Row=17; E=7; F=$null},
    @{Row=18; E=10; F=20251208},
    @{Row=19; E=10; F=20251208},
    @{Row=20; E=10; F=20251208},
    @{Row=21; E=10; F=20251208},
    @{Row=22; E=7; F=$null},
    @{Row=23; E=7; F=$null},
    @{Row=24; E=7; F=$null},
    @{Row=25; E=7; F=$null},
    @{Row=26; E=7; F=$null},
    @{Row=27; E=1; F=$null},
    @{Row=28; E=10; F=20251208},
    @{Row=29; E=10; F=20251208},
    @{Row=30; E=10; F=20251208},
    @{Row=31; E=10; F=20251208},
    @{Row=32; E=10; F=20251208},
    @{Row=33; E=10; F=20251208},
    @{Row=34; E=10; F=20251208},
    @{Row=35; E=10; F=20251208},
    @{Row=37; E=10; F=20251208},
    @{Row=38; E=10; F=20251208},
    @{Row=39; E=10; F=20251208},
    @{Row=40; E=7; F=20251208},
    @{Row=41; E=7; F=20251208},
    @{Row=42; E=10; F=20251208},
    @{Row=43; E=7; F=$null},
    @{Row=44; E=7; F=20251208},
    @{Row=45; E=7; F=$null},
    @{Row=46; E=7; F=20251208},
    @{Row=47; E=10; F=20251208},
    @{Row=48; E=7; F=20251208},
    @{Row=49; E=1; F=$null},
    @{Row=50; E=5; F=$null},
    @{Row=51; E=5; F=$null},
    @{Row=52; E=5; F=$null},
    @{Row=53; E=5; F=$null},
    @{Row=54; E=5; F=$null},
    @{Row=55; E=5; F=$null},
    @{Row=56; E=5; F=$null},
    @{Row=57; E=5; F=$null},
    @{Row=58; E=9; F=$null},
    @{Row=59; E=9; F=$null},
    @{Row=60; E=9; F=$null},
    @{Row=61; E=1; F=$null},
    @{Row=62; E=9; F=$null},
    @{Row=63; E=9; F=$null},
    @{Row=64; E=9; F=$null},
    @{Row=65; E=10; F=20251208},
    @{Row=66; E=10; F=20251208},
    @{Row=67; E=10; F=20251208},
    @{Row=68; E=10; F=20251208},
    @{Row=69; E=10; F=20251208},
    @{Row=70; E=1; F=$null},
    @{Row=71; E=1; F=$null},
    @{Row=72; E=1; F=$null},
    @{Row=73; E=1; F=$null},
    @{Row=74; E=1; F=$null},
    @{Row=75; E=1; F=$null},
    @{Row=76; E=1; F=$null},
    @{Row=77; E=4; F=$null},
    @{Row=78; E=4; F=$null},
    @{Row=79; E=4; F=$null},
    @{Row=80; E=4; F=$null},
    @{Row=81; E=4; F=$null},
    @{Row=82; E=4; F=$null},
    @{Row=83; E=4; F=$null},
    @{Row=84; E=4; F=$null},
    @{Row=85; E=4; F=$null},
    @{Row=86; E=4; F=$null},
    @{Row=87; E=7; F=20251208},
    @{Row=88; E=7; F=20251208},
    @{Row=89; E=7; F=20251208},
    @{Row=90; E=7; F=20251208},
    @{Row=91; E=7; F=$null},
    @{Row=92; E=7; F=20251208},
    @{Row=93; E=4; F=$null},
    @{Row=94; E=3; F=$null},
    @{Row=95; E=3; F=$null},
    @{Row=96; E=1; F=$null},
    @{Row=97; E=1; F=$null},
    @{Row=98; E=1; F=$null},
    @{Row=99; E=1; F=$null}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    if ($u.F -ne $null) {
        $ws.Cells.Item($u.Row, 6).Value = $u.F
    }
}
